# Update "Förändrad" (changed) date column C for rows 2-10 from
# serial date 45224 (2023-10-25) to 45233 (2023-11-03).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45224) {
        $cell.Value2 = 45233
    }
}
